# Generate Report for Handoff
# "b.md" has now been handed off for localization (zh-cn and de-de):
#  - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#  - New handoff xliff files were generated for b.md (zh-cn + de-de)
#  - Content Duplicate flag flips from True to False for b.md
#  - An error is now reported: handback file version is stale vs. the new b.md handoff
#  - The "Error Detail" column is widened to fit the new message

$wb = $excel.ActiveWorkbook

$newStatus      = "Ready for handoff"
$newHoDate      = "2016-09-06 02:41:40"

$zhHandoffFile  = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate  = "2016-09-06 02:41:36"

$deHandoffFile  = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate  = "2016-09-06 02:41:40"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/21c5198ef5c1bc9e4a0956fa570dc88280f711fd/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b4e681decb0950d2b520f57d1553fa1a736da3f5/e2e/b.md."

# ---- Overview sheet: row 3 is "b.md" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value2 = $newStatus
$wsOverview.Range("F3").Value2 = $newStatus
$wsOverview.Range("G3").Value2 = $newHoDate

# ---- zh-cn sheet: row 3 is "b.md" ----
# NB: "False" is entered with a leading apostrophe so it stays a literal
# text value ("Content Duplicate" is a text column in this sheet, not a
# real boolean column) instead of Excel's automatic Boolean auto-detect.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value2 = $newStatus
$wsZh.Range("F3").Value2 = "'False"
$wsZh.Range("G3").Value2 = $zhHandoffFile
$wsZh.Range("H3").Value2 = $zhHandoffDate
$wsZh.Range("P3").Value2 = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet: row 3 is "b.md" ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value2 = $newStatus
$wsDe.Range("F3").Value2 = "'False"
$wsDe.Range("G3").Value2 = $deHandoffFile
$wsDe.Range("H3").Value2 = $deHandoffDate
$wsDe.Range("P3").Value2 = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664
